$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row data: row number, date string (with dashes), D, E, F, G, H
$rows = @(
    @{R=3;  Date="28-07-2022"; D=1; E=0; F=0; G=1; H=1},
    @{R=4;  Date="01-08-2022"; D=1; E=1; F=0; G=0; H=0},
    @{R=5;  Date="04-08-2022"; D=1; E=1; F=0; G=0; H=0},
    @{R=6;  Date="08-08-2022"; D=0; E=0; F=0; G=0; H=1},
    @{R=7;  Date="11-08-2022"; D=0; E=0; F=0; G=0; H=1},
    @{R=8;  Date="15-08-2022"; D=0; E=0; F=0; G=0; H=1},
    @{R=9;  Date="18-08-2022"; D=0; E=0; F=0; G=0; H=1},
    @{R=10; Date="22-08-2022"; D=0; E=0; F=0; G=0; H=1},
    @{R=11; Date="25-08-2022"; D=1; E=1; F=0; G=0; H=0},
    @{R=12; Date="29-08-2022"; D=1; E=1; F=0; G=0; H=0},
    @{R=13; Date="01-09-2022"; D=1; E=1; F=0; G=0; H=0},
    @{R=14; Date="05-09-2022"; D=0; E=0; F=0; G=0; H=1},
    @{R=15; Date="08-09-2022"; D=1; E=1; F=0; G=0; H=0},
    @{R=16; Date="12-09-2022"; D=0; E=0; F=0; G=0; H=1},
    @{R=17; Date="15-09-2022"; D=0; E=0; F=0; G=0; H=1},
    @{R=18; Date="19-09-2022"; D=0; E=0; F=0; G=0; H=1},
    @{R=19; Date="22-09-2022"; D=0; E=0; F=0; G=0; H=1},
    @{R=20; Date="26-09-2022"; D=0; E=0; F=0; G=0; H=1},
    @{R=21; Date="29-09-2022"; D=0; E=0; F=0; G=0; H=1}
)

foreach ($row in $rows) {
    $r = $row.R
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row.Date
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
}

$wb.Save()
